# codigo_sap_clientes.xlsx — "Añadido configuracion multiple por proveedor,
# nueva configuracion para ESP, formateo de codigo"
#
# Concrete, COM-reachable content changes from the diff:
#   1. Sheet1!B2: correct/reformat the SAP code 7001353 -> 70001353
#      (adds the missing leading "0" so it matches the 8-digit pattern
#      used by every other row in the sap_code column).
#   2. The sheet's active/selected cell moves from D9 to B15.
#
# (The diff also touches pure window/session metadata emitted by Excel on
# save - the x15ac:absPath OneDrive/SharePoint URL, the xr:revisionPtr
# save-revision GUID, and the workbookView x/y/width/height screen
# position. These are host-chrome bookkeeping values, not workbook
# content, and are out of reach of the Excel object model exposed here -
# they are window-manager/session state, not Range/Cell/Worksheet state.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1) Reformat the SAP code in B2 so it carries the leading zero like the
#    rest of the column (7001353 -> 70001353).
$ws.Range("B2").Value = 70001353

# 2) Move the active selection to B15 (was D9).
$ws.Activate()
$ws.Range("B15").Select()
